$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so numeric-looking strings (e.g. "408.00")
# are not silently coerced into numbers and lose formatting/precision.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.554.84'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.386.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.00'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.48'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +9.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.596'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.675'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.121'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.81'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.55%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.908.02'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.43'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.78'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.384.53'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.516.51'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.05'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '85.18'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '315.76'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.83'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +11.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.37'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '29.60'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.67'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.173'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.55'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.38'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.68%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.86'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0482'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.95'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.44'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.94'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '138.60'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.98'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.294'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.03'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.79'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.58%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.52'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.131.23'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.38%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.97'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.50%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.29'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.82%  '
